$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the price series. Insert a row before
# the current row 185 (shifting all subsequent rows, including the last
# one, down by one) and populate it with the new observation.
$ws.Rows("185").Insert()

$ws.Range("A185").Value = 10
$ws.Range("B185").Value = "Vega Modelo de Temuco"
$ws.Range("C185").Value = "La Araucanía"
$ws.Range("D185").Value = 44813
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = 100112039
$ws.Range("G185").Value = "Ciboulette"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 30
$ws.Range("K185").Value = 8000
$ws.Range("L185").Value = 8000
$ws.Range("M185").Value = 8000
$ws.Range("N185").Value = "`$/docena de atados"
$ws.Range("O185").Value = "Provincia de Cautín"
$ws.Range("P185").Value = 2667
$ws.Range("Q185").Value = 3
$ws.Range("R185").Value = "Hortaliza"
